$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("'2023-01-26", "21:00 ", "Atletico de Madrid", "1–3", "Real Madrid", "César Soto", 1883, 320),
    @("'2023-01-29", "16:15 ", "Atletico de Madrid", "1–0", "CA Osasuna", "Alberola Rojas", 624, 317),
    @("'2023-02-04", "18:30 ", "Atletico de Madrid", "1–1", "Getafe CF", "Antonio Matéu", 2830, 317)
)

$startRow = 30
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $startRow + $i
    $rowData = $data[$i]

    $ws.Cells.Item($row, 1).Value = $rowData[0]
    $ws.Cells.Item($row, 1).Style = "Normal"

    $ws.Cells.Item($row, 2).Value = $rowData[1]
    $ws.Cells.Item($row, 3).Value = $rowData[2]
    $ws.Cells.Item($row, 4).Value = $rowData[3]
    $ws.Cells.Item($row, 5).Value = $rowData[4]
    $ws.Cells.Item($row, 6).Value = $rowData[5]
    $ws.Cells.Item($row, 7).Value = $rowData[6]
    $ws.Cells.Item($row, 8).Value = $rowData[7]
}
